# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.06328177979961902, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569)
    3 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
    4 = @(0.7287194209349384, 1.65323645889881, 16.98373111632243, 0.4998867070740569)
    5 = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569)
    6 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]
    $sum = $b + $c + $d + $e

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = $sum
}
